$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 803, shifting all
# subsequent rows (803-836) down by one (to 804-837) and extending the
# sheet's used range from A1:T836 to A1:T837.
$ws.Rows.Item(803).Insert()

# Populate the newly inserted row 803 with the new weekly record.
$ws.Range("A803").Value = 9
$ws.Range("B803").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C803").Value = "Metropolitana"
$ws.Range("D803").Value = 44747
$ws.Range("E803").Value = 13
$ws.Range("F803").Value = "Fruta"
$ws.Range("G803").Value = 100104
$ws.Range("H803").Value = "Frutos de pepita"
$ws.Range("I803").Value = 100104005
$ws.Range("J803").Value = "Pera"
$ws.Range("K803").Value = "Packham's Triumph"
$ws.Range("L803").Value = "Primera"
$ws.Range("M803").Value = 380
$ws.Range("N803").Value = 10000
$ws.Range("O803").Value = 10000
$ws.Range("P803").Value = 10000
$ws.Range("Q803").Value = "$/caja 18 kilos granel"
$ws.Range("R803").Value = "Paine"
$ws.Range("S803").Value = 556
$ws.Range("T803").Value = 18
